$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value, append Source/Target rows ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B8").Value = "2025-05-02T13:48:14+00:00"

# Copy formatting from the last existing data row (14) down into the two
# new rows so the new cells pick up style index "2" like their neighbours.
$meta.Range("A14:B14").Copy()
$meta.Range("A15:B15").PasteSpecial(-4122)
$meta.Range("A16:B16").PasteSpecial(-4122)

$meta.Range("A15").Value = "Source"
$meta.Range("B15").Value = "https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/SubmissionSet"
$meta.Range("A16").Value = "Target"
$meta.Range("B16").Value = "https://interop.esante.gouv.fr/ig/fhir/pdsm/StructureDefinition/pdsm-submissionset-comprehensive"
